# Apply the "simulator full-month coverage, persist logs, fix employees" edit.
#
# 1. Fix the employee/client names on the "Weekly Timesheet" sheet (rows 2-6).
# 2. Fix the Employee ID on the "Jason Schema" sheet (rows 2-6).
# 3. Populate the previously-zeroed Rate/Total columns on both sheets for the
#    simulated days (now that the simulator runs the full month and persists
#    its logs, rate/total are no longer 0).

$wb = $excel.ActiveWorkbook

# ---- Sheet 1: "Weekly Timesheet" ----
$ws1 = $wb.Worksheets.Item("Weekly Timesheet")

# Client / employee name corrections (column B, rows 2-6)
$ws1.Range("B2").Value = "Tubergen"
$ws1.Range("B3").Value = "Hewett"
$ws1.Range("B4").Value = "Durfee"
$ws1.Range("B5").Value = "Markfield"
$ws1.Range("B6").Value = "Corr"

# Rate (E) / Total (F) for each daily row, previously 0
$ws1.Range("E2").Value = 150
$ws1.Range("F2").Value = 1200
$ws1.Range("E3").Value = 150
$ws1.Range("F3").Value = 1200
$ws1.Range("E4").Value = 150
$ws1.Range("F4").Value = 1200
$ws1.Range("E5").Value = 150
$ws1.Range("F5").Value = 1200
$ws1.Range("E6").Value = 150
$ws1.Range("F6").Value = 1200

# Subtotal / grand total rows
$ws1.Range("F8").Value = 6000
$ws1.Range("F12").Value = 6000
$ws1.Range("F13").Value = 6000

# ---- Sheet 2: "Jason Schema" ----
$ws2 = $wb.Worksheets.Item("Jason Schema")

# Employee ID correction (column B, rows 2-6)
$ws2.Range("B2").Value = "emp_35u1tnme"
$ws2.Range("B3").Value = "emp_35u1tnme"
$ws2.Range("B4").Value = "emp_35u1tnme"
$ws2.Range("B5").Value = "emp_35u1tnme"
$ws2.Range("B6").Value = "emp_35u1tnme"

# Client name corrections (column D, rows 2-6) -- mirrors Weekly Timesheet col B
$ws2.Range("D2").Value = "Tubergen"
$ws2.Range("D3").Value = "Hewett"
$ws2.Range("D4").Value = "Durfee"
$ws2.Range("D5").Value = "Markfield"
$ws2.Range("D6").Value = "Corr"

$ws2.Range("F2").Value = 150
$ws2.Range("G2").Value = 1200
$ws2.Range("F3").Value = 150
$ws2.Range("G3").Value = 1200
$ws2.Range("F4").Value = 150
$ws2.Range("G4").Value = 1200
$ws2.Range("F5").Value = 150
$ws2.Range("G5").Value = 1200
$ws2.Range("F6").Value = 150
$ws2.Range("G6").Value = 1200
